$d = $word.ActiveDocument

# 1) "Salvare gli appunti preferiti nel proprio account di Google Docs, creare nuovi appunti"
#    -> "Creare nuovi appunti"
$d.Content.Find.Execute("Salvare gli appunti preferiti nel proprio account di Google Docs, creare nuovi appunti", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Creare nuovi appunti", 2) | Out-Null

# 2) Insert new sentence about cartolerie after "modificare i propri appunti caricati."
#    and before " Inviare e-mail agli sviluppatori..."
$d.Content.Find.Execute("modificare i propri appunti caricati. Inviare", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "modificare i propri appunti caricati. Visualizzare sulla mappa le cartolerie in prossimità del luogo di studio. Inviare", 2) | Out-Null

# 3) " Login locale" -> "Login locale" (drop leading space)
$d.Content.Find.Execute(" Login locale", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Login locale", 2) | Out-Null

# 4) user data bullet: drop ", l'avatar che lo rappresenta" and "rating" -> "level"
$d.Content.Find.Execute("il luogo in cui studia, la sua biografia, l’avatar che lo rappresenta e il trust rating nella community.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "il luogo in cui studia, la sua biografia e il trust level nella community.", 2) | Out-Null

# 5) remove the empty paragraph right after "...star rating." and before the "Funzionalità" heading
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "`r") {
        $prevText = $p.Previous().Range.Text
        if ($prevText -like "*star rating.`r") {
            $target = $p
        }
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 6) rewrite the "; in più possono ..." sentence
$d.Content.Find.Execute("; in più possono decidere di trasferire i loro appunti preferiti sulla piattaforma Google Docs per averli sempre a portata di mano e possono contattare gli sviluppatori per fornire preziosi feedback o richiedere assistenza. Chi si registra può modificare il trust level degli altri utenti per premiare gli utenti che condividono gli appunti migliori. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "; in più possono visualizzare le cartolerie nei pressi del loro luogo di studio e possono contattare gli sviluppatori per fornire preziosi feedback o richiedere assistenza. Il trust level di ogni utente viene modificato in base ad un algoritmo del sito che valuta l’utente in base alla qualità dei suoi appunti. ", 2) | Out-Null
